$d = $word.ActiveDocument

# Resources menu: the "RPC Explorer" link label was renamed to
# "Insight Explorer" (Chinese Simplified translation pass).
$d.Content.Find.Execute(
    "RPC Explorer", $true, $false, $false, $false, $false,
    $true, 1, $false, "Insight Explorer", 2
) | Out-Null
